$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) and Volume(1h) (E) columns for the rows we touch
# so numeric-looking strings (e.g. "69.632.32", "0.0890", "1.00") are preserved
# exactly as text instead of being auto-coerced into numbers by Excel.
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D6").NumberFormat = "@"
$ws.Range("D10:D11").NumberFormat = "@"
$ws.Range("D13:D16").NumberFormat = "@"
$ws.Range("D18:D19").NumberFormat = "@"
$ws.Range("D21:D24").NumberFormat = "@"
$ws.Range("D26:D29").NumberFormat = "@"
$ws.Range("D31:D37").NumberFormat = "@"
$ws.Range("D39:D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E2:E37").NumberFormat = "@"
$ws.Range("E39:E51").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = '69.632.32'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").Value = '3.670.74'
$ws.Range("E3").Value = '  -0.79%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '623.75'
$ws.Range("E5").Value = '  -7.33%  '
$ws.Range("D6").Value = '160.02'
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("E9").Value = '  -2.13%  '
$ws.Range("D10").Value = '7.15'
$ws.Range("E10").Value = '  +0.67%  '
$ws.Range("D11").Value = '0.439'
$ws.Range("E11").Value = '  -1.19%  '
$ws.Range("E12").Value = '  -2.55%  '
$ws.Range("D13").Value = '4.291.11'
$ws.Range("E13").Value = '  -0.73%  '
$ws.Range("D14").Value = '32.48'
$ws.Range("E14").Value = '  -1.04%  '
$ws.Range("D15").Value = '3.733.10'
$ws.Range("E15").Value = '  +1.05%  '
$ws.Range("D16").Value = '69.680.33'
$ws.Range("E16").Value = '  +0.01%  '
$ws.Range("E17").Value = '  +0.63%  '
$ws.Range("D18").Value = '6.51'
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("D19").Value = '15.86'
$ws.Range("E19").Value = '  -2.10%  '
$ws.Range("E20").Value = '  +5.50%  '
$ws.Range("D21").Value = '470.46'
$ws.Range("E21").Value = '  -0.85%  '
$ws.Range("D22").Value = '0.649'
$ws.Range("E22").Value = '  -0.72%  '
$ws.Range("D23").Value = '79.90'
$ws.Range("E23").Value = '  -0.70%  '
$ws.Range("D24").Value = '3.818.16'
$ws.Range("E24").Value = '  -0.76%  '
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").Value = '0.0000124'
$ws.Range("E26").Value = '  -2.03%  '
$ws.Range("D27").Value = '11.04'
$ws.Range("E27").Value = '  +0.27%  '
$ws.Range("D28").Value = '8.69'
$ws.Range("E28").Value = '  -5.03%  '
$ws.Range("D29").Value = '2.58'
$ws.Range("E29").Value = '  -3.96%  '
$ws.Range("E30").Value = '  -4.44%  '
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("D32").Value = '1.98'
$ws.Range("E32").Value = '  -1.97%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '26.60'
$ws.Range("E33").Value = '  -1.09%  '
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = '0.164'
$ws.Range("E34").Value = '  -0.80%  '
$ws.Range("D35").Value = '6.39'
$ws.Range("E35").Value = '  -3.19%  '
$ws.Range("D36").Value = '3.674.05'
$ws.Range("E36").Value = '  -0.41%  '
$ws.Range("D37").Value = '8.30'
$ws.Range("E37").Value = '  -2.81%  '
$ws.Range("D39").Value = '178.35'
$ws.Range("E39").Value = '  +3.29%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '5.82'
$ws.Range("E41").Value = '  -5.14%  '
$ws.Range("D42").Value = '2.20'
$ws.Range("E42").Value = '  -1.53%  '
$ws.Range("D43").Value = '0.0890'
$ws.Range("E43").Value = '  -1.86%  '
$ws.Range("D44").Value = '0.925'
$ws.Range("E44").Value = '  -1.64%  '
$ws.Range("D45").Value = '46.74'
$ws.Range("E45").Value = '  -0.65%  '
$ws.Range("D46").Value = '29.01'
$ws.Range("E46").Value = '  +3.66%  '
$ws.Range("D47").Value = '2.72'
$ws.Range("E47").Value = '  -2.01%  '
$ws.Range("D48").Value = '7.87'
$ws.Range("E48").Value = '  -0.29%  '
$ws.Range("D49").Value = '0.000264'
$ws.Range("E49").Value = '  -5.47%  '
$ws.Range("E50").Value = '  -5.94%  '
$ws.Range("D51").Value = '1.21'
$ws.Range("E51").Value = '  -5.66%  '
